# Apply updated price/profit figures to the Leve tables across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1755.1666
$ws.Range("I40").Value = 1249.4445
$ws.Range("J40").Value = 2260.889
$ws.Range("K40").Value = 1249.4445
$ws.Range("L40").Value = 2260.889
$ws.Range("M40").Value = -1074.4445
$ws.Range("N40").Value = -2610.889
$ws.Range("H69").Value = 13890.111
$ws.Range("I69").Value = 7003.6665
$ws.Range("K69").Value = 21010.9995
$ws.Range("M69").Value = -20136.9995
$ws.Range("H72").Value = 13890.111
$ws.Range("I72").Value = 7003.6665
$ws.Range("K72").Value = 63032.9985
$ws.Range("M72").Value = -58664.9985
$ws.Range("H92").Value = 750.4706
$ws.Range("I92").Value = 773.6
$ws.Range("J92").Value = 717.4286
$ws.Range("K92").Value = 773.6
$ws.Range("L92").Value = 717.4286
$ws.Range("M92").Value = 474.4
$ws.Range("N92").Value = -3213.4286
$ws.Range("H100").Value = 3077.7778
$ws.Range("I100").Value = 3000.1428
$ws.Range("K100").Value = 3000.1428
$ws.Range("M100").Value = -2459.1428
$ws.Range("H103").Value = 100000490
$ws.Range("I103").Value = 100000490
$ws.Range("K103").Value = 300001470
$ws.Range("M103").Value = -300000884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2074
$ws.Range("I2").Value = 2106.25
$ws.Range("K2").Value = 2106.25
$ws.Range("M2").Value = -1993.25
$ws.Range("H32").Value = 11194.826
$ws.Range("I32").Value = 11570.954
$ws.Range("K32").Value = 11570.954
$ws.Range("M32").Value = -11283.954
$ws.Range("H61").Value = 1499.909
$ws.Range("I61").Value = 1499.909
$ws.Range("K61").Value = 1499.909
$ws.Range("M61").Value = -1287.909
$ws.Range("H74").Value = 1661
$ws.Range("J74").Value = 1832.6666
$ws.Range("L74").Value = 1832.6666
$ws.Range("N74").Value = -3580.6666
$ws.Range("H77").Value = 1661
$ws.Range("J77").Value = 1832.6666
$ws.Range("L77").Value = 9163.333000000001
$ws.Range("N77").Value = -17899.333
$ws.Range("H116").Value = 2074
$ws.Range("I116").Value = 2106.25
$ws.Range("K116").Value = 2106.25
$ws.Range("M116").Value = 187.75
$ws.Range("H122").Value = 2263.5217
$ws.Range("I122").Value = 2252.5908
$ws.Range("K122").Value = 6757.7724
$ws.Range("M122").Value = -4307.7724
$ws.Range("H132").Value = 3617.8462
$ws.Range("I132").Value = 1673.3334
$ws.Range("J132").Value = 5284.5713
$ws.Range("K132").Value = 5020.0002
$ws.Range("L132").Value = 15853.7139
$ws.Range("M132").Value = -2490.0002
$ws.Range("N132").Value = -20913.7139
$ws.Range("H136").Value = 1499.909
$ws.Range("I136").Value = 1499.909
$ws.Range("K136").Value = 4499.727000000001
$ws.Range("M136").Value = -1949.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2074
$ws.Range("I3").Value = 2106.25
$ws.Range("K3").Value = 2106.25
$ws.Range("M3").Value = -1992.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 222.77777
$ws.Range("J7").Value = 97.5
$ws.Range("L7").Value = 97.5
$ws.Range("N7").Value = -323.5
$ws.Range("H29").Value = 2999.5
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2707
$ws.Range("H31").Value = 2668.75
$ws.Range("I31").Value = 1908.5
$ws.Range("J31").Value = 4949.5
$ws.Range("K31").Value = 1908.5
$ws.Range("L31").Value = 4949.5
$ws.Range("M31").Value = -1613.5
$ws.Range("N31").Value = -5539.5
$ws.Range("H33").Value = 2032
$ws.Range("I33").Value = 1665
$ws.Range("J33").Value = 3500
$ws.Range("K33").Value = 1665
$ws.Range("L33").Value = 3500
$ws.Range("M33").Value = -1286
$ws.Range("N33").Value = -4258
$ws.Range("H34").Value = 2668.75
$ws.Range("I34").Value = 1908.5
$ws.Range("J34").Value = 4949.5
$ws.Range("K34").Value = 1908.5
$ws.Range("L34").Value = 4949.5
$ws.Range("M34").Value = -1706.5
$ws.Range("N34").Value = -5353.5
$ws.Range("H62").Value = 3042.8333
$ws.Range("J62").Value = 3119.5
$ws.Range("L62").Value = 3119.5
$ws.Range("N62").Value = -4367.5
$ws.Range("H65").Value = 3042.8333
$ws.Range("J65").Value = 3119.5
$ws.Range("L65").Value = 15597.5
$ws.Range("N65").Value = -21837.5
$ws.Range("H99").Value = 2886.4
$ws.Range("I99").Value = 2873.7778
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2873.7778
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1375.7778
$ws.Range("N99").Value = -5996
$ws.Range("H105").Value = 1639.4
$ws.Range("I105").Value = 799.25
$ws.Range("K105").Value = 799.25
$ws.Range("M105").Value = 947.75
$ws.Range("H126").Value = 2886.4
$ws.Range("I126").Value = 2873.7778
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8621.3334
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6151.3334
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.64706
$ws.Range("I2").Value = 63.375
$ws.Range("K2").Value = 380.25
$ws.Range("M2").Value = -267.25
$ws.Range("H33").Value = 175
$ws.Range("J33").Value = 100
$ws.Range("L33").Value = 600
$ws.Range("N33").Value = -1166
$ws.Range("H40").Value = 53.46154
$ws.Range("I40").Value = 34.6
$ws.Range("J40").Value = 116.333336
$ws.Range("K40").Value = 138.4
$ws.Range("L40").Value = 465.333344
$ws.Range("M40").Value = -69.40000000000001
$ws.Range("N40").Value = -603.333344
$ws.Range("H86").Value = 1113.75
$ws.Range("I86").Value = 999.6667
$ws.Range("K86").Value = 2999.0001
$ws.Range("M86").Value = -1813.0001
$ws.Range("H89").Value = 1113.75
$ws.Range("I89").Value = 999.6667
$ws.Range("K89").Value = 8997.0003
$ws.Range("M89").Value = -3069.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4135.4165
$ws.Range("I7").Value = 4247.273
$ws.Range("K7").Value = 4247.273
$ws.Range("M7").Value = -4135.273
$ws.Range("H22").Value = 2599.9
$ws.Range("I22").Value = 2175
$ws.Range("K22").Value = 2175
$ws.Range("M22").Value = -1880
$ws.Range("H27").Value = 2599.9
$ws.Range("I27").Value = 2175
$ws.Range("K27").Value = 2175
$ws.Range("M27").Value = -2068
$ws.Range("H55").Value = 1399.8572
$ws.Range("I55").Value = 1249.8
$ws.Range("K55").Value = 1249.8
$ws.Range("M55").Value = -1076.8
$ws.Range("H122").Value = 4259.6
$ws.Range("I122").Value = 4266.5
$ws.Range("K122").Value = 12799.5
$ws.Range("M122").Value = -10349.5
$ws.Range("H126").Value = 4135.4165
$ws.Range("I126").Value = 4247.273
$ws.Range("K126").Value = 12741.819
$ws.Range("M126").Value = -10271.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H74").Value = 24998
$ws.Range("J74").Value = 24998
$ws.Range("L74").Value = 24998
$ws.Range("N74").Value = -26870
$ws.Range("H77").Value = 24998
$ws.Range("J77").Value = 24998
$ws.Range("L77").Value = 74994
$ws.Range("N77").Value = -84354
$ws.Range("H81").Value = 2197.8572
$ws.Range("I81").Value = 2197.8572
$ws.Range("K81").Value = 4395.7144
$ws.Range("M81").Value = -3334.7144
$ws.Range("H84").Value = 2197.8572
$ws.Range("I84").Value = 2197.8572
$ws.Range("K84").Value = 21978.572
$ws.Range("M84").Value = -16674.572
